# NewCards.xlsx — "Added PickCard Button, Changed persistent buff logic"
#
# The only logical (data-level) change in the target workbook is a new
# card row appended to the "Cards" sheet: a Rare, Stamina-cost card named
# "再一次" ("Again") whose Effect1 is 0 (VDoublePlayEffect — "the next
# card played is used twice"), mirroring the existing "网络波动" card but
# as its own purchasable/unique entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cards")

# New row 19 — ID 17 "再一次" ("Again")
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "再一次 "
$ws.Range("C19").Value = "下一张使用的牌会使用两次 "
$ws.Range("D19").Value = "Rare"
$ws.Range("F19").Value = "Stamina"
$ws.Range("H19").Value = 6
$ws.Range("I19").Value = 5
$ws.Range("J19").Value = 1
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
